$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets
# ---------------------------------------------------------------------
$wsReviews = $wb.Worksheets.Item(1)
$wsHistory = $wb.Worksheets.Item(2)

# Rename the reviews sheet (LH-TC-REGISTERATION-Reviews -> LH-TC-USERHOME-Reviews)
$wsReviews.Name = "LH-TC-USERHOME-Reviews"

# ---------------------------------------------------------------------
# Reviews sheet: close out the existing review rows (Owner Status /
# Reviewer verification both flip to "Closed") and append the new
# review row for LH_TC_USERHOME_REVIEW_004.
# ---------------------------------------------------------------------
$wsReviews.Range("I2").Value = "Closed"
$wsReviews.Range("J2").Value = "Closed"

$wsReviews.Range("I3").Value = "Closed"
$wsReviews.Range("J3").Value = "Closed"

$wsReviews.Range("I4").Value = "Closed"
$wsReviews.Range("J4").Value = "Closed"
$wsReviews.Range("J4").VerticalAlignment = -4108

$wsReviews.Range("B5").Value = "LH_TC_USERHOME_REVIEW_004"
$wsReviews.Range("C5").Value = "LH_TC_USERHOME.xlsx file"
$wsReviews.Range("D5").Value = "Ahmed Abuzaid"
$wsReviews.Range("E5").Value = "v2.1"
$wsReviews.Range("F5").Value = 'ia all data you mention that there is a correct password but this password doesn''t match acceptance criteria of the password'
$wsReviews.Range("G5").Value = 'so please make the password valid something like that "CorrectPassword@123"'
$wsReviews.Range("H5").Value = "hala"
$wsReviews.Range("I5").Value = "Open"
$wsReviews.Range("J5").Value = "Open"

$wsReviews.Range("C5").HorizontalAlignment = -4108
$wsReviews.Range("C5").VerticalAlignment = -4108
$wsReviews.Range("J5").HorizontalAlignment = -4108
$wsReviews.Range("J5").VerticalAlignment = -4108

$wsReviews.Rows.Item(5).RowHeight = 105

# ---------------------------------------------------------------------
# Version History sheet: fix casing of the v1.0 entry, close the
# v1.1 entry and record who closed it / when / with what comment.
# ---------------------------------------------------------------------
$wsHistory.Range("C2").Value = "review  the user home feature test cases"

$wsHistory.Rows.Item(3).Delete()
$wsHistory.Rows.Item(3).Insert()
$wsHistory.Range("A3").Value = "v1.1"
$wsHistory.Range("B3").Value = "Ahmed Abuzaid"
$wsHistory.Range("C3").Value = "close the previous status and add more comments "
$wsHistory.Range("D3").Value = "14/5/2025"
$wsHistory.Rows.Item(3).RowHeight = 30

$loHistory = $wsHistory.ListObjects.Item(1)
$loHistory.Resize($wsHistory.Range("A1:D3"))

# ---------------------------------------------------------------------
# Window / selection state: Version History becomes the active tab,
# each sheet keeps its own last-used selection.
# ---------------------------------------------------------------------
$wsReviews.Range("G7").Select()
$wsHistory.Activate()
$wsHistory.Range("C13").Select()

Write-Host "done"
